$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 10 with a new timesheet entry (matching the date/time
# formatting already used by the other rows in columns A-C)
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("B10").NumberFormat = "h:mm"
$ws.Range("C10").NumberFormat = "h:mm"

$ws.Range("A10").Value = 42646
$ws.Range("B10").Value = 0.90277777777777779
$ws.Range("C10").Value = 0.95138888888888884
$ws.Range("D10").Value = "Implementação do pagamento por Créditos"

# Move the active selection to D11, as in the committed sheet view
$ws.Range("D11").Select()

$wb.Save()
